# Update column G ("K") values on Sheet1 for rows 2-71.
# These were regenerated from source data ("use K instead of Strike#,
# regen std/mean, calc and write s_vals") and are simply replaced with
# their newly computed literal values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @(2,5,0,0,1,1,0,1,1,2,0,1,2,1,1,2,0,3,0,2,2,1,0,1,1,0,0,0,1,1,5,2,2,2,1,1,2,2,0,1,1,2,0,0,0,0,1,1,1,0,3,0,2,1,3,0,0,1,0,2,1,1,3,0,2,1,1,0,1,2)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newValues[$i]
}
